$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - use same style as other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F8 - time_taken values, no special style (like other data cells)
$ws.Range("F2").Value = "2021-10-05 13:39:36.296221"
$ws.Range("F3").Value = "2021-10-05 13:39:36.296231"
$ws.Range("F4").Value = "2021-10-05 13:39:36.296234"
$ws.Range("F5").Value = "2021-10-05 13:39:36.296237"
$ws.Range("F6").Value = "2021-10-05 13:39:36.296240"
$ws.Range("F7").Value = "2021-10-05 13:39:36.296242"
$ws.Range("F8").Value = "2021-10-05 13:39:36.296245"
